$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-ThinBox($rng) {
    $rng.Borders.Color = 0
    $rng.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
    $rng.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
}

# --- Row "3" style cells first (plain, xf6) ---
$a3 = $ws.Range("A70")
Set-ThinBox $a3
$d3 = $ws.Range("D70")
Set-ThinBox $d3
$e3 = $ws.Range("E70")
Set-ThinBox $e3

# --- Header row A2:F2 bold (xf7) ---
$hdr = $ws.Range("A71:F71")
Set-ThinBox $hdr
$hdr.Font.Bold = $true

# --- Number format cells B3,C3 (xf8) ---
$b3 = $ws.Range("B70")
Set-ThinBox $b3
$b3.NumberFormat = "#,##0"
$c3 = $ws.Range("C70")
Set-ThinBox $c3
$c3.NumberFormat = "#,##0"

# --- Text format F3 (xf9) ---
$f3 = $ws.Range("F70")
Set-ThinBox $f3
$f3.NumberFormat = "@"

# --- Merged header row A1:F1 bold+centered (xf10 initial, then split into border2/3/4) ---
$full1 = $ws.Range("A72:F72")
Set-ThinBox $full1
$full1.Font.Bold = $true
$full1.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$mid1 = $ws.Range("B72:E72")
$mid1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlInsideVertical).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$mid1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$mid1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$a1 = $ws.Range("A72")
$a1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$f1 = $ws.Range("F72")
$f1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
